# 1st commit - Data-driven - Creating a common data provider
#
# Adds a new "OpenAccountTest" worksheet after the existing
# "AddCustomerTest" sheet, makes it the active/selected sheet, and
# populates it with a small customer/currency data table.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Insert the new sheet right after AddCustomerTest so it becomes sheet 2
# (sheetId 2, rId2) and is the active tab.
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "OpenAccountTest"

# Header row (bold, like the existing AddCustomerTest headers)
$ws2.Range("A1").Value = "customer"
$ws2.Range("B1").Value = "currency"
$ws2.Range("A1:B1").Font.Bold = $true

# Data row
$ws2.Range("A2").Value = "Janielle Joy Gregorio"
$ws2.Range("B2").Value = "Dollar"

# Column A sized to fit the customer name
$ws2.Columns.Item(1).ColumnWidth = 18.6666666667

# Selections: AddCustomerTest loses its previous selection/tab focus,
# OpenAccountTest becomes the active sheet with B3 selected.
$ws1.Range("A2").Select() | Out-Null
$ws2.Range("B3").Select() | Out-Null
$ws2.Activate() | Out-Null
